# Optuna Attempt (go back with original)
#
# Restores a previous round of forecast numbers on the "Forecast Comparison"
# sheet (MyForecast / Inventory Coverage / Seasonality Index for weeks
# W8-W23) and the dependent roll-up figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---------------------------------------------
# Row 2 (W8) - Seasonality Index only
$wsForecast.Range("L2").Value = 0.82

# Row 3 (W9) - Seasonality Index only
$wsForecast.Range("L3").Value = 1.04

# Row 4 (W10) - Seasonality Index only
$wsForecast.Range("L4").Value = 0.9

# Row 5 (W11) - Seasonality Index only
$wsForecast.Range("L5").Value = 0.88

# Row 6 (W12)
$wsForecast.Range("D6").Value = 74
$wsForecast.Range("H6").Value = 18.07
$wsForecast.Range("L6").Value = 1.13

# Row 7 (W13)
$wsForecast.Range("D7").Value = 72
$wsForecast.Range("H7").Value = 17.49
$wsForecast.Range("L7").Value = 0.84

# Row 8 (W14)
$wsForecast.Range("D8").Value = 71
$wsForecast.Range("H8").Value = 16.76
$wsForecast.Range("L8").Value = 0.95

# Row 9 (W15)
$wsForecast.Range("D9").Value = 69
$wsForecast.Range("H9").Value = 16.21
$wsForecast.Range("L9").Value = 0.99

# Row 10 (W16)
$wsForecast.Range("D10").Value = 68
$wsForecast.Range("H10").Value = 15.53
$wsForecast.Range("L10").Value = 0.87

# Row 11 (W17)
$wsForecast.Range("D11").Value = 68
$wsForecast.Range("H11").Value = 14.4
$wsForecast.Range("L11").Value = 0.8100000000000001

# Row 12 (W18)
$wsForecast.Range("D12").Value = 67
$wsForecast.Range("H12").Value = 13.76
$wsForecast.Range("L12").Value = 1.01

# Row 13 (W19)
$wsForecast.Range("D13").Value = 65
$wsForecast.Range("H13").Value = 13.15
$wsForecast.Range("L13").Value = 1.12

# Row 14 (W20)
$wsForecast.Range("D14").Value = 76
$wsForecast.Range("H14").Value = 10.35
$wsForecast.Range("L14").Value = 1.11

# Row 15 (W21)
$wsForecast.Range("D15").Value = 63
$wsForecast.Range("H15").Value = 11.25
$wsForecast.Range("L15").Value = 0.91

# Row 16 (W22)
$wsForecast.Range("D16").Value = 74
$wsForecast.Range("H16").Value = 8.699999999999999
$wsForecast.Range("L16").Value = 1.17

# Row 17 (W23)
$wsForecast.Range("D17").Value = 74
$wsForecast.Range("H17").Value = 7.73
$wsForecast.Range("L17").Value = 0.85

# --- Summary sheet -----------------------------------------------------------
# These cells hold numeric-looking values stored as TEXT (e.g. "1136"), as in
# the source file. Assigning a plain numeric-looking string to a cell whose
# number format is "General" makes Excel coerce it into a real number, so the
# format is switched to Text ("@") just long enough to force string storage,
# then cleared again so the cell keeps its original (default/General) style -
# matching the source file, which never set an explicit number format on
# these cells.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1136"
$wsSummary.Range("B9").ClearFormats()

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "578"
$wsSummary.Range("B10").ClearFormats()

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "76"
$wsSummary.Range("B12").ClearFormats()

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "64"
$wsSummary.Range("B14").ClearFormats()
